$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Add([System.Type]::Missing, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = "Feuil2"
$ws2.Range("A1").Value = "feuil2 A1"
$ws2.Range("A2").Value = "feuil2 A2"
$ws2.Range("A4").Value = "feuil2 A4"
$ws2.Range("B2").Value = "feuil2 B2"
$ws2.Range("B3").Value = "feuil2 B3"
$ws2.Range("B4").Value = "feuil2 B4"
$ws2.Range("B1").Value = "feuil2 B1"
$ws2.Range("A1:B4").SetPhonetic()
$ws2.Range("D4").Select()
